$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 102
$ws.Range("B3").Value = "WoLunZengYa"
$ws.Range("A4").Value = 103
$ws.Range("B4").Value = "JiyouLengQue"
$ws.Range("A5").Value = 104
$ws.Range("B5").Value = "JieQiMenTi"
$ws.Range("A6").Value = 105
$ws.Range("B6").Value = "JinQiQiGuan"
$ws.Range("A7").Value = 106
$ws.Range("B7").Value = "RanYouPenShe"
$ws.Range("A8").Value = 107
$ws.Range("B8").Value = "ChuShuiKou"
$ws.Range("A9").Value = 108
$ws.Range("B9").Value = "ZhenKongBeng"
$ws.Range("A10").Value = 109
$ws.Range("B10").Value = "DianHuoXianQuan"
$ws.Range("A11").Value = 110
$ws.Range("B11").Value = "RanYouDaoGui"
$ws.Range("A12").Value = 201
$ws.Range("B12").Value = "TuLunZhouGai"
$ws.Range("A13").Value = 202
$ws.Range("B13").Value = "PeiQiZhengShi"
$ws.Range("A14").Value = 203
$ws.Range("B14").Value = "NiuZhuanJianZhen"
$ws.Range("A15").Value = 204
$ws.Range("B15").Value = "FaDongJiQianGai"
$ws.Range("A16").Value = 205
$ws.Range("B16").Value = "ZhengShiLian"
$ws.Range("A17").Value = 206
$ws.Range("B17").Value = "TuLunZhouZhengShiChiLun"
$ws.Range("A18").Value = 301
$ws.Range("B18").Value = "JinQiTuLunZhou"
$ws.Range("A19").Value = 302
$ws.Range("B19").Value = "PaiQiTuLunZhou"
$ws.Range("A20").Value = 303
$ws.Range("B20").Value = "QiMenYaoBi"
$ws.Range("A21").Value = 304
$ws.Range("B21").Value = "JianXiTiaoJieQi"
$ws.Range("A22").Value = 401
$ws.Range("B22").Value = "QiGangGai"
$ws.Range("A23").Value = 501
$ws.Range("B23").Value = "YouDiKe"
$ws.Range("A24").Value = 502
$ws.Range("B24").Value = "JiYouBengChuanDongLianTiao"
$ws.Range("A25").Value = 503
$ws.Range("B25").Value = "JiYouBeng"
$ws.Range("A26").Value = 504
$ws.Range("B26").Value = "JiYouDaoLiuBan"
$ws.Range("A27").Value = 601
$ws.Range("B27").Value = "LianGanZhouChengGai"
$ws.Range("A28").Value = 602
$ws.Range("B28").Value = "HuoSaiLianGan"
$ws.Range("A29").Value = 603
$ws.Range("B29").Value = "QiHuan"
$ws.Range("A30").Value = 604
$ws.Range("B30").Value = "YouHuan"
$ws.Range("A31").Value = 701
$ws.Range("B31").Value = "XiaBuQuZhouXiang"
$ws.Range("A32").Value = 702
$ws.Range("B32").Value = "QuZhou"

$ws.Range("B32").Select() | Out-Null
